$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Week Of" banner text (row 3)
$ws.Range("A3").Value = "Week Of:  April 22 - 26, 2019"

# Update the action item labels (only the two that changed text; "Game Logic" and
# "Configuring game scenes and adding addition features" stay the same)
$ws.Range("A12").Value = "GUI elements and Enchancing certain elements"
$ws.Range("A11").Value = "Networking features and connection between server and Clients"

# Row 11 (Networking features...) - Last Week / This Week / Next Week
$ws.Range("D11").Value = "Learn Node JS"
$ws.Range("E11").Value = "Start Server for Game"

# Row 10 (Configuring game scenes...) - Last Week / This Week / Next Week
$ws.Range("D10").Value = "Learn Node JS and HTML"
$ws.Range("E10").Value = "Decide how Tank will function in game"

# Row 12 (GUI elements...) - This Week / Last Week
$ws.Range("E12").Value = "Create ideas for how the GUI and user interface will look like"
$ws.Range("D12").Value = "Learn Node JS, HTML, and CSS"

# Row 10 / Row 12 - Next Week
$ws.Range("F10").Value = "Implement the Tank in the game"
$ws.Range("F12").Value = "Implement GUI and UI deisgns"

# Row 11 - Next Week
$ws.Range("F11").Value = "Configure the Server and Client Logic together and make them connect"

# Row 9 (Game Logic) - This Week / Next Week / Last Week
$ws.Range("E9").Value = "Start deciding how game will be played out"
$ws.Range("F9").Value = "Write Game Logic Code"
$ws.Range("D9").Value = "Learn Node JS and Brainstorm Game functions"

# Remove the now-unused "Issues" column entries for rows 9-12
$ws.Range("G9").Clear()
$ws.Range("G10").Clear()
$ws.Range("G11").Clear()
$ws.Range("G12").Clear()

# Update the selected cell / scroll position to match the saved view state
$ws.Range("D10").Select()
